$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024483449217156
$ws.Range("D2").Value = 1.030238759434109
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.035560546575582
$ws.Range("I2").Value = 1.03340948704682
$ws.Range("J2").Value = 1.029657809676376
$ws.Range("K2").Value = 1.033050348551973
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.038356800559069
$ws.Range("N2").Value = 1.013943662199267

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025223686306697
$ws.Range("D3").Value = 1.030789887109413
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.036576525314912
$ws.Range("I3").Value = 1.033561745334133
$ws.Range("J3").Value = 1.030038055000321
$ws.Range("K3").Value = 1.033410480799456
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.039181650541821
$ws.Range("N3").Value = 1.014069792001033

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025703142704302
$ws.Range("D4").Value = 1.031146877991444
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.03723483652632
$ws.Range("I4").Value = 1.033659320998349
$ws.Range("J4").Value = 1.030283885486492
$ws.Range("K4").Value = 1.033643178011274
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.039715682454944
$ws.Range("N4").Value = 1.014151320401119

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02590481781283
$ws.Range("D5").Value = 1.03129704486568
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.037511805934682
$ws.Range("I5").Value = 1.033700114773508
$ws.Range("J5").Value = 1.030387180364309
$ws.Range("K5").Value = 1.033740923169737
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.039940259795977
$ws.Range("N5").Value = 1.014185573956615

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025938686447578
$ws.Range("D6").Value = 1.03132226364918
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.037558322909351
$ws.Range("I6").Value = 1.033706950895639
$ws.Range("J6").Value = 1.030404520932154
$ws.Range("K6").Value = 1.033757330242893
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.039977971408432
$ws.Range("N6").Value = 1.014191324037406

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02570583706046
$ws.Range("D7").Value = 1.031148884185331
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.037238536560441
$ws.Range("I7").Value = 1.033659866979665
$ws.Range("J7").Value = 1.030285265924433
$ws.Range("K7").Value = 1.033644484405616
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.039718682993551
$ws.Range("N7").Value = 1.014151778181748

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024733516732176
$ws.Range("D8").Value = 1.030424937046335
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.035903713176993
$ws.Range("I8").Value = 1.033461138856267
$ws.Range("J8").Value = 1.029786358939329
$ws.Range("K8").Value = 1.033172125064131
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.038635499869876
$ws.Range("N8").Value = 1.013986305931158

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.023023866542729
$ws.Range("D9").Value = 1.029152193716395
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.033558572421373
$ws.Range("I9").Value = 1.033103740836073
$ws.Range("J9").Value = 1.028905633574785
$ws.Range("K9").Value = 1.032337272936467
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.0367291350725
$ws.Range("N9").Value = 1.013694081196959

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021886691469261
$ws.Range("D10").Value = 1.028305780710677
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.031999920134912
$ws.Range("I10").Value = 1.032860664119079
$ws.Range("J10").Value = 1.028317481876686
$ws.Range("K10").Value = 1.031779096965868
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.035459871667156
$ws.Range("N10").Value = 1.013498856709165

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021394918405619
$ws.Range("D11").Value = 1.027939790305005
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.031326154791976
$ws.Range("I11").Value = 1.032754276166458
$ws.Range("J11").Value = 1.028062582051324
$ws.Range("K11").Value = 1.031537033868337
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.034910671991443
$ws.Range("N11").Value = 1.01341423048936

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021212348568475
$ws.Range("D12").Value = 1.027803923866964
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.031076061059598
$ws.Range("I12").Value = 1.032714589194782
$ws.Range("J12").Value = 1.027967868027398
$ws.Range("K12").Value = 1.031447066559496
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.034706736002615
$ws.Range("N12").Value = 1.013382783006051

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02125150602018
$ws.Range("D13").Value = 1.02783306408967
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.031129699178525
$ws.Range("I13").Value = 1.032723109856129
$ws.Range("J13").Value = 1.027988185977037
$ws.Range("K13").Value = 1.031466367290836
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.03475047817523
$ws.Range("N13").Value = 1.013389529204242

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021379825145782
$ws.Range("D14").Value = 1.027928557927291
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.031305478435532
$ws.Range("I14").Value = 1.032750999088117
$ws.Range("J14").Value = 1.028054753624853
$ws.Range("K14").Value = 1.031529598243369
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.034893813323285
$ws.Range("N14").Value = 1.013411631304459

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021458899643208
$ws.Range("D15").Value = 1.027987405312984
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.031413804760746
$ws.Range("I15").Value = 1.032768160090196
$ws.Range("J15").Value = 1.028095763826137
$ws.Range("K15").Value = 1.031568549760394
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.034982134986349
$ws.Range("N15").Value = 1.01342524735582

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021919342275805
$ws.Range("D16").Value = 1.028330081227259
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.032044659874685
$ws.Range("I16").Value = 1.03286770089699
$ws.Range("J16").Value = 1.028334394066652
$ws.Range("K16").Value = 1.031795154217695
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.035496328762613
$ws.Range("N16").Value = 1.013504471152637

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.022208336003505
$ws.Range("D17").Value = 1.028545171106109
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.032440685071724
$ws.Range("I17").Value = 1.032929837028587
$ws.Range("J17").Value = 1.02848402072238
$ws.Range("K17").Value = 1.031937199186049
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.035818976780468
$ws.Range("N17").Value = 1.013554141599705

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022376961945198
$ws.Range("D18").Value = 1.028670678676594
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.032671789912647
$ws.Range("I18").Value = 1.03296597051272
$ws.Range("J18").Value = 1.028571273513861
$ws.Range("K18").Value = 1.032020015909929
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.036007210402636
$ws.Range("N18").Value = 1.013583104557791

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022434469291775
$ws.Range("D19").Value = 1.028713481820075
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.032750609318934
$ws.Range("I19").Value = 1.032978272509511
$ws.Range("J19").Value = 1.028601020711531
$ws.Range("K19").Value = 1.032048248181008
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.036071399696957
$ws.Range("N19").Value = 1.013592978643997

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.022177323411234
$ws.Range("D20").Value = 1.028522088887332
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.032398183924745
$ws.Range("I20").Value = 1.032923181728107
$ws.Range("J20").Value = 1.028467969456006
$ws.Range("K20").Value = 1.031921962789693
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.035784355711809
$ws.Range("N20").Value = 1.01354881335612

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02134203568167
$ws.Range("D21").Value = 1.027900435185688
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.031253711047785
$ws.Range("I21").Value = 1.032742791085171
$ws.Range("J21").Value = 1.028035152010419
$ws.Range("K21").Value = 1.031510979793264
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.034851602995258
$ws.Range("N21").Value = 1.013405123157359

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020817416746343
$ws.Range("D22").Value = 1.027510033004365
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.030535135419856
$ws.Range("I22").Value = 1.0326283905303
$ws.Range("J22").Value = 1.027762832824555
$ws.Range("K22").Value = 1.031252264908757
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.03426549914409
$ws.Range("N22").Value = 1.01331470129841

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021095473545097
$ws.Range("D23").Value = 1.027716948757317
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.030915970763677
$ws.Range("I23").Value = 1.032689129277674
$ws.Range("J23").Value = 1.027907211985262
$ws.Range("K23").Value = 1.031389443894414
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.03457616986354
$ws.Range("N23").Value = 1.013362642910867

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.022191336479148
$ws.Range("D24").Value = 1.028532518596391
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.032417388023712
$ws.Range("I24").Value = 1.032926189310242
$ws.Range("J24").Value = 1.028475222400185
$ws.Range("K24").Value = 1.031928847571828
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.035799999363658
$ws.Range("N24").Value = 1.01355122098782

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023465402787058
$ws.Range("D25").Value = 1.029480868848186
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.034164011508593
$ws.Range("I25").Value = 1.03319698766603
$ws.Range("J25").Value = 1.029133503527663
$ws.Range("K25").Value = 1.032553390696885
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.037221691130707
$ws.Range("N25").Value = 1.013769701911783
